$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": append a new date column AQ (26-jul) mirroring the
# existing AP (25-jul) column - same header style + the day's hourly prices.
# ---------------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the header cell's formatting (bold/border/centered) from AP1 onto AQ1,
# then set the new header text.
$wsSpot.Range("AP1").Copy()
$wsSpot.Range("AQ1").PasteSpecial(-4122)
$wsSpot.Range("AQ1").Value = "26-jul"

$spotValues = @{
    2  = 96.27
    3  = 77.16
    4  = 66.28
    5  = 56.84
    6  = 46.4
    7  = 50.74
    8  = 57.33
    9  = 55.01
    10 = 41.05
    11 = 33.2
    12 = 17.69
    13 = 13.78
    14 = 18
    15 = 14.58
    16 = 8.67
    17 = 8.15
    18 = 5.98
    19 = 17.3
    20 = 27.13
    21 = 49.33
    22 = 65.29000000000001
    23 = 92.56999999999999
    24 = 100
    25 = 90.66
}

foreach ($row in $spotValues.Keys) {
    # AQ is column 43
    $wsSpot.Cells.Item($row, 43).Value = $spotValues[$row]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append row 40 with the 2025-07-24 price. Force the date cell
# to stay plain text (matching the existing A2:A39 cells) instead of letting
# it auto-convert to a date serial, then drop the temporary text format so
# no stray style sticks to the cell.
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$gazDate = $wsGaz.Cells.Item(40, 1)
$gazDate.NumberFormat = "@"
$gazDate.Value = "2025-07-24"
$gazDate.ClearFormats()
$wsGaz.Cells.Item(40, 2).Value = 32.075

# ---------------------------------------------------------------------------
# Sheet "CO2": append row 40 with the 2025-07-24 price (same treatment).
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$co2Date = $wsCo2.Cells.Item(40, 1)
$co2Date.NumberFormat = "@"
$co2Date.Value = "2025-07-24"
$co2Date.ClearFormats()
$wsCo2.Cells.Item(40, 2).Value = 70.2
